$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.540.65"
$ws.Range("E2").Value = "  -0.89%  "

$ws.Range("D3").Value = "2.342.91"
$ws.Range("E3").Value = "  -1.40%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.66"
$ws.Range("E5").Value = "  -0.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.25"
$ws.Range("E6").Value = "  -3.39%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.517"
$ws.Range("E8").Value = "  -3.07%  "

$ws.Range("D9").Value = "2.339.59"
$ws.Range("E9").Value = "  -1.29%  "

$ws.Range("E10").Value = "  -0.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.157"
$ws.Range("E11").Value = "  +1.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.31"
$ws.Range("E12").Value = "  +1.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.343"
$ws.Range("E13").Value = "  +0.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.49"
$ws.Range("E14").Value = "  -3.78%  "

$ws.Range("D15").Value = "2.748.88"
$ws.Range("E15").Value = "  -1.84%  "

$ws.Range("D16").Value = "60.461.91"
$ws.Range("E16").Value = "  -0.69%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000160"
$ws.Range("E17").Value = "  -3.63%  "

$ws.Range("D18").Value = "2.345.85"
$ws.Range("E18").Value = "  -1.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.61"
$ws.Range("E19").Value = "  -1.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "320.60"
$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.54"
$ws.Range("E22").Value = "  -2.71%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.36"
$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.64"
$ws.Range("E25").Value = "  -14.79%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.33"
$ws.Range("E26").Value = "  +7.25%  "

$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").Value = "2.466.45"
$ws.Range("E28").Value = "  -0.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.00"
$ws.Range("E29").Value = "  -0.46%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.147"
$ws.Range("E30").Value = "  +1.47%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.36"
$ws.Range("E31").Value = "  -5.01%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "488.06"
$ws.Range("E32").Value = "  -8.71%  "

$ws.Range("D33").Value = "0.0₃0852"
$ws.Range("E33").Value = "  -9.98%  "

$ws.Range("E34").Value = "  -1.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.49"
$ws.Range("E35").Value = "  -4.67%  "

$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.60"
$ws.Range("E37").Value = "  -0.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.376"
$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.41"
$ws.Range("E39").Value = "  +2.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.83"
$ws.Range("E40").Value = "  +4.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.23"
$ws.Range("E41").Value = "  -5.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "144.31"
$ws.Range("E42").Value = "  +4.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "141.87"
$ws.Range("E44").Value = "  +0.97%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.55"
$ws.Range("E45").Value = "  -0.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.02"
$ws.Range("E46").Value = "  -11.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0515"
$ws.Range("E47").Value = "  -0.52%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.566"
$ws.Range("E48").Value = "  -1.55%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.96"
$ws.Range("E49").Value = "  -6.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0900"
$ws.Range("E50").Value = "  -1.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0220"
$ws.Range("E51").Value = "  -2.08%  "
